# New crime data collected - update weekly CompStat figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume/issue number and reporting week dates --------------
$ws.Range("A8").Value = "Volume 30   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/25/2023  Through  12/31/2023"

# --- Row 16 (Robbery) ---------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 133.333333333333
$ws.Range("I16").Value = 102
$ws.Range("J16").Value = 116
$ws.Range("K16").Value = -12.068965517241
$ws.Range("L16").Value = -11.304347826087
$ws.Range("M16").Value = 13.333333333333
$ws.Range("N16").Value = -83.278688524590

# --- Row 17 (Fel. Assault) ----------------------------------------------------
# D17/E17 flip from the text placeholders ("0" / "***.*") to real numbers, so
# pull the numeric-cell formatting from a sibling cell that already has it.
$ws.Range("D17").Value = 2
$ws.Range("C17").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").Value = 50
$ws.Range("K17").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Range("F17").Value = 11
$ws.Range("H17").Value = 175
$ws.Range("I17").Value = 110
$ws.Range("J17").Value = 102
$ws.Range("K17").Value = 7.843137254901
$ws.Range("L17").Value = 27.906976744186
$ws.Range("M17").Value = 92.982456140350
$ws.Range("N17").Value = 11.111111111111

# --- Row 18 (Burglary) --------------------------------------------------------
$ws.Range("C18").Value = 4
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -11.111111111111
$ws.Range("I18").Value = 111
$ws.Range("J18").Value = 136
$ws.Range("K18").Value = -18.382352941176
$ws.Range("L18").Value = 37.037037037037
$ws.Range("M18").Value = 9.900990099009
$ws.Range("N18").Value = -87.947882736156

# --- Row 19 (Gr. Larceny) -----------------------------------------------------
$ws.Range("C19").Value = 9
$ws.Range("E19").Value = -43.75
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = 3.846153846153
$ws.Range("I19").Value = 744
$ws.Range("J19").Value = 820
$ws.Range("K19").Value = -9.268292682926
$ws.Range("L19").Value = -0.268096514745
$ws.Range("M19").Value = 8.931185944363
$ws.Range("N19").Value = -58.620689655172

# --- Row 20 (G.L.A.) -----------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 5
$ws.Range("C20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = -60
$ws.Range("K20").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -42.857142857142
$ws.Range("I20").Value = 94
$ws.Range("J20").Value = 75
$ws.Range("K20").Value = 25.333333333333
$ws.Range("L20").Value = 8.045977011494
$ws.Range("M20").Value = 184.848484848485
$ws.Range("N20").Value = -91.847354726799

# --- Row 21 (TOTAL) -------------------------------------------------------------
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = 13.75
$ws.Range("I21").Value = 1175
$ws.Range("J21").Value = 1263
$ws.Range("K21").Value = -6.967537608867
$ws.Range("L21").Value = 4.723707664884
$ws.Range("M21").Value = 20.389344262295
$ws.Range("N21").Value = -74.467622772707

# --- Row 22 (Transit) -----------------------------------------------------------
# C22 flips from a numeric 1 to the text placeholder "0" - borrow that look
# (value + format) from another cell in the sheet that already carries it.
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4163) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("J22").Value = 28
$ws.Range("K22").Value = 0

# --- Row 23 (Housing) ------------------------------------------------------------
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4163) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null

$ws.Range("D23").Value = 1
$ws.Range("C17").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = -100
$ws.Range("K17").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null

$ws.Range("G23").Value = 1
$ws.Range("C17").Copy() | Out-Null
$ws.Range("G23").PasteSpecial(-4122) | Out-Null
$ws.Range("H23").Value = 700
$ws.Range("K17").Copy() | Out-Null
$ws.Range("H23").PasteSpecial(-4122) | Out-Null

$ws.Range("J23").Value = 25
$ws.Range("K23").Value = 72
$ws.Range("L23").Value = 59.259259259259
$ws.Range("M23").Value = 59.259259259259

# --- Row 24 (Petit Larceny) -------------------------------------------------------
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -23.529411764705
$ws.Range("F24").Value = 87
$ws.Range("G24").Value = 67
$ws.Range("H24").Value = 29.850746268656
$ws.Range("I24").Value = 1171
$ws.Range("J24").Value = 1222
$ws.Range("K24").Value = -4.173486088379
$ws.Range("L24").Value = -6.767515923566
$ws.Range("M24").Value = 12.487992315081

# --- Row 25 (Misd. Assault) --------------------------------------------------------
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = -31.25
$ws.Range("I25").Value = 222
$ws.Range("J25").Value = 223
$ws.Range("K25").Value = -0.448430493273
$ws.Range("L25").Value = 17.460317460317
$ws.Range("M25").Value = -19.272727272727

# --- Row 27 (Other Sex Crimes) ------------------------------------------------------
$ws.Range("C27").Value = 2
$ws.Range("C17").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").Value = 1
$ws.Range("C17").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = 100
$ws.Range("K17").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 43
$ws.Range("J27").Value = 54
$ws.Range("K27").Value = -20.370370370370
$ws.Range("L27").Value = -23.214285714285

# --- Row 28 (Shooting Vic.) -----------------------------------------------------------
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4163) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null

# --- Row 29 (Shooting Inc.) -----------------------------------------------------------
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4163) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null
